# The Pearson logo picture (embedded in footer1.xml / footer2.xml) is
# currently named "image1.png" on both its wp:docPr and pic:cNvPr
# elements; it should be renamed to "image2.png".
#
# The BTec logo picture (embedded in header1.xml / header2.xml) is
# currently named "image2.jpg" on both its wp:docPr and pic:cNvPr
# elements; it should be renamed to "image1.jpg".
#
# Neither the real media part filenames (media/image1.png,
# media/image2.jpg) nor any relationship ids change - only the cosmetic
# "name" attribute that Word stores on the drawing's docPr/cNvPr nodes.
#
# These name attributes aren't reachable through the regular
# InlineShape object (InlineShape has no writable Name), so we go
# through the document's flattened WordOpenXML, patch the two name
# values there, and write the whole thing back.

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

$xml = $xml.Replace('name="image1.png"', 'name="image2.png"')
$xml = $xml.Replace('name="image2.jpg"', 'name="image1.jpg"')

$d.WordOpenXML = $xml
